$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we touch so numeric-looking
# strings (e.g. "183.50", "0.999") are preserved verbatim as text,
# matching the original inlineStr cell type instead of being
# auto-coerced to a Number by Excel (which would drop formatting
# like trailing zeros).
$cells = @('D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'D6', 'E6', 'D7', 'E7', 'D8', 'E8', 'E9', 'D10', 'E10', 'B11', 'C11', 'D11', 'E11', 'B12', 'C12', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'E17', 'D18', 'E18', 'B19', 'C19', 'D19', 'E19', 'B20', 'C20', 'D20', 'E20', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'D28', 'E28', 'D29', 'E29', 'D30', 'E30', 'D31', 'E31', 'E32', 'E33', 'D34', 'E34', 'B35', 'C35', 'D35', 'E35', 'B36', 'C36', 'D36', 'E36', 'D37', 'E37', 'E38', 'B39', 'C39', 'D39', 'E39', 'D40', 'E40', 'B41', 'C41', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'E45', 'D46', 'E46', 'E47', 'B48', 'C48', 'D48', 'E48', 'B49', 'C49', 'D49', 'E49', 'B50', 'C50', 'D50', 'E50', 'B51', 'C51', 'D51', 'E51')
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '65.838.47'
$ws.Range('E2').Value = '  -5.67%  '
$ws.Range('D3').Value = '3.540.01'
$ws.Range('E3').Value = '  -6.06%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '581.87'
$ws.Range('E5').Value = '  -6.19%  '
$ws.Range('D6').Value = '183.50'
$ws.Range('E6').Value = '  +1.06%  '
$ws.Range('D7').Value = '3.531.06'
$ws.Range('E7').Value = '  -6.25%  '
$ws.Range('D8').Value = '0.608'
$ws.Range('E8').Value = '  -4.42%  '
$ws.Range('E9').Value = '  +0.50%  '
$ws.Range('D10').Value = '0.659'
$ws.Range('E10').Value = '  -8.93%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = '0.143'
$ws.Range('E11').Value = '  -12.93%  '
$ws.Range('B12').Value = 'Avalanche'
$ws.Range('C12').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D12').Value = '53.16'
$ws.Range('E12').Value = '  -8.18%  '
$ws.Range('D13').Value = '0.0000249'
$ws.Range('E13').Value = '  -16.24%  '
$ws.Range('D14').Value = '9.65'
$ws.Range('E14').Value = '  -10.76%  '
$ws.Range('D15').Value = '4.109.63'
$ws.Range('E15').Value = '  -5.78%  '
$ws.Range('D16').Value = '3.542.14'
$ws.Range('E16').Value = '  -6.10%  '
$ws.Range('E17').Value = '  -1.01%  '
$ws.Range('D18').Value = '18.11'
$ws.Range('E18').Value = '  -8.04%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '12.00'
$ws.Range('E19').Value = '  -8.21%  '
$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D20').Value = '65.588.08'
$ws.Range('E20').Value = '  -5.73%  '
$ws.Range('E21').Value = '  -8.93%  '
$ws.Range('D22').Value = '391.75'
$ws.Range('E22').Value = '  -6.20%  '
$ws.Range('D23').Value = '4.28'
$ws.Range('E23').Value = '  -9.72%  '
$ws.Range('D24').Value = '84.47'
$ws.Range('E24').Value = '  -6.09%  '
$ws.Range('D25').Value = '2.85'
$ws.Range('E25').Value = '  -7.17%  '
$ws.Range('D26').Value = '12.30'
$ws.Range('E26').Value = '  -4.11%  '
$ws.Range('D27').Value = '6.02'
$ws.Range('E27').Value = '  -1.04%  '
$ws.Range('D28').Value = '10.32'
$ws.Range('E28').Value = '  -7.23%  '
$ws.Range('D29').Value = '3.51'
$ws.Range('E29').Value = '  -8.73%  '
$ws.Range('D30').Value = '8.83'
$ws.Range('E30').Value = '  -9.12%  '
$ws.Range('D31').Value = '30.63'
$ws.Range('E31').Value = '  -8.36%  '
$ws.Range('E32').Value = '  -8.35%  '
$ws.Range('E33').Value = '  -6.25%  '
$ws.Range('D34').Value = '607.95'
$ws.Range('E34').Value = '  -2.68%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.110'
$ws.Range('E35').Value = '  -8.41%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').Value = '61.83'
$ws.Range('E36').Value = '  -7.81%  '
$ws.Range('D37').Value = '40.49'
$ws.Range('E37').Value = '  -9.61%  '
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').Value = '0.368'
$ws.Range('E39').Value = '  -9.29%  '
$ws.Range('D40').Value = '0.0₃0732'
$ws.Range('E40').Value = '  -18.13%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = '0.997'
$ws.Range('E41').Value = '  -0.30%  '
$ws.Range('D42').Value = '0.128'
$ws.Range('E42').Value = '  -9.71%  '
$ws.Range('D43').Value = '2.877.77'
$ws.Range('E43').Value = '  +0.99%  '
$ws.Range('D44').Value = '2.74'
$ws.Range('E44').Value = '  -11.80%  '
$ws.Range('E45').Value = '  -9.01%  '
$ws.Range('D46').Value = '2.41'
$ws.Range('E46').Value = '  -10.00%  '
$ws.Range('E47').Value = '  -6.19%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').Value = '3.03'
$ws.Range('E48').Value = '  -4.70%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').Value = '138.51'
$ws.Range('E49').Value = '  -2.64%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = '8.34'
$ws.Range('E50').Value = '  -10.81%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = '2.71'
$ws.Range('E51').Value = '  -3.08%  '
